$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "02.29"
$ws.Name = "02.29"

# Put the header value in A2 (A1 stays empty, used range starts at A2)
$ws.Range("A2").Value = "時刻"

# Select A2 to match the saved selection state
$ws.Range("A2").Select()
